# Iran and canada predictions updated
#
# The "Canada" block (prediction day 2021-01-09, previously rows 46-52)
# gets a brand new sibling block inserted right above it for prediction
# day 2021-01-02 (new rows 46-52); the old block shifts down to rows
# 53-59 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing "2021-01-09" Canada rows (46-52) down to 53-59,
# opening up a fresh block of 7 blank rows at 46-52.
$ws.Rows("46:52").Insert()

# Fill in the new "2021-01-02" Canada prediction block.
$ws.Range("A46").Value = "2021-01-02"
$ws.Range("B46").Value = "03 Jan -- 09 Jan 2021"
$ws.Range("C46").Value = 94.56999999999999
$ws.Range("D46").Value = 226.44
$ws.Range("E46").Value = 131.87
$ws.Range("F46").Value = "KNN"
$ws.Range("G46").Value = 1.23
$ws.Range("H46").Value = 69.72
$ws.Range("I46").Value = 86.05
$ws.Range("J46").Value = 130.59
$ws.Range("K46").Value = 128.92

$ws.Range("A47").Value = "2021-01-02"
$ws.Range("B47").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D47").Value = 246.14
$ws.Range("F47").Value = "KNN"

$ws.Range("A48").Value = "2021-01-02"
$ws.Range("B48").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D48").Value = 260.58
$ws.Range("F48").Value = "KNN"

$ws.Range("A49").Value = "2021-01-02"
$ws.Range("B49").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D49").Value = 262.97
$ws.Range("F49").Value = "KNN"

$ws.Range("A50").Value = "2021-01-02"
$ws.Range("B50").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D50").Value = 254.31
$ws.Range("F50").Value = "KNN"

$ws.Range("A51").Value = "2021-01-02"
$ws.Range("B51").Value = "07 Feb -- 13 Feb 2021"
$ws.Range("D51").Value = 247.89
$ws.Range("F51").Value = "KNN"

$ws.Range("A52").Value = "2021-01-02"
$ws.Range("B52").Value = "14 Feb -- 20 Feb 2021"
$ws.Range("D52").Value = 230.15
$ws.Range("F52").Value = "KNN"

$ws.Dimension = "A1:K59"
